# Lab 1 - Team Formation: reposition the GitHub Classroom link textbox on
# the last slide (sldId 262) slightly to the left.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$sh.Left = 2348582 / 12700
